# Scale the "value" column (D) from units of 10,000 to absolute units
# (multiply by 10000) for all data rows.
#
# A plain floating-point multiplication (`$x * 10000`) can land on a
# different (but numerically adjacent) double than the one you'd get by
# parsing the decimal-shifted text directly, because the source values
# carry 6 decimal digits and binary doubles cannot represent them
# exactly. To reproduce the exact target doubles, the scaling is done by
# round-tripping through .NET's `decimal` (base-10, exact) arithmetic on
# the clean/shortest string form of each value, then parsing that exact
# decimal text back into a double.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 33; $row++) {
    $cell = $ws.Cells.Item($row, 4)
    $orig = $cell.Value2
    $s = $orig.ToString()
    $d = [decimal]$s
    $scaled = $d * 10000
    $cell.Value2 = [double]$scaled.ToString()
}
